$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update odds values (columns C, D, E) for jornada update
$ws.Range("D1").Value = 3.5
$ws.Range("E1").Value = 4

$ws.Range("C2").Value = 2.0499999999999998
$ws.Range("D2").Value = 3.4
$ws.Range("E2").Value = 3.4

$ws.Range("C4").Value = 5.5
$ws.Range("D4").Value = 3.75
$ws.Range("E4").Value = 1.61

$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 3.4
$ws.Range("E5").Value = 3.5

$ws.Range("C6").Value = 1.87
$ws.Range("D6").Value = 3.5
$ws.Range("E6").Value = 4

$ws.Range("C7").Value = 2.9
$ws.Range("D7").Value = 3.4
$ws.Range("E7").Value = 2.2999999999999998

$ws.Range("C8").Value = 1.71
$ws.Range("D8").Value = 3.6
$ws.Range("E8").Value = 4.75

$ws.Range("D9").Value = 3.4
$ws.Range("E9").Value = 2.9

$ws.Range("C10").Value = 1.87
$ws.Range("D10").Value = 3.4

# Update the view: scroll back to top-left A1 and move selection to E10
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E10").Select()
